$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-05-06 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-05-07 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("38×30=", $true, $false, $false, $false, $false, $true, 1, $false, "69×39=", 2) | Out-Null
$d.Content.Find.Execute("69×69=", $true, $false, $false, $false, $false, $true, 1, $false, "88×90=", 2) | Out-Null
$d.Content.Find.Execute("82×18=", $true, $false, $false, $false, $false, $true, 1, $false, "64×14=", 2) | Out-Null
$d.Content.Find.Execute("40×87=", $true, $false, $false, $false, $false, $true, 1, $false, "84×22=", 2) | Out-Null
$d.Content.Find.Execute("27×87=", $true, $false, $false, $false, $false, $true, 1, $false, "23×66=", 2) | Out-Null
$d.Content.Find.Execute("89×53=", $true, $false, $false, $false, $false, $true, 1, $false, "97×20=", 2) | Out-Null
$d.Content.Find.Execute("20×15=", $true, $false, $false, $false, $false, $true, 1, $false, "93×97=", 2) | Out-Null
$d.Content.Find.Execute("21×55=", $true, $false, $false, $false, $false, $true, 1, $false, "25×73=", 2) | Out-Null
$d.Content.Find.Execute("64×62=", $true, $false, $false, $false, $false, $true, 1, $false, "47×55=", 2) | Out-Null
$d.Content.Find.Execute("84×79=", $true, $false, $false, $false, $false, $true, 1, $false, "32×30=", 2) | Out-Null
$d.Content.Find.Execute("55×55=", $true, $false, $false, $false, $false, $true, 1, $false, "47×50=", 2) | Out-Null
$d.Content.Find.Execute("52×96=", $true, $false, $false, $false, $false, $true, 1, $false, "85×12=", 2) | Out-Null
$d.Content.Find.Execute("24×73=", $true, $false, $false, $false, $false, $true, 1, $false, "22×14=", 2) | Out-Null
$d.Content.Find.Execute("32×67=", $true, $false, $false, $false, $false, $true, 1, $false, "44×55=", 2) | Out-Null
$d.Content.Find.Execute("54×34=", $true, $false, $false, $false, $false, $true, 1, $false, "27×18=", 2) | Out-Null
$d.Content.Find.Execute("73×43=", $true, $false, $false, $false, $false, $true, 1, $false, "22×49=", 2) | Out-Null
$d.Content.Find.Execute("51×65=", $true, $false, $false, $false, $false, $true, 1, $false, "39×41=", 2) | Out-Null
$d.Content.Find.Execute("68×23=", $true, $false, $false, $false, $false, $true, 1, $false, "80×13=", 2) | Out-Null
$d.Content.Find.Execute("88×98=", $true, $false, $false, $false, $false, $true, 1, $false, "82×99=", 2) | Out-Null
$d.Content.Find.Execute("87×67=", $true, $false, $false, $false, $false, $true, 1, $false, "57×68=", 2) | Out-Null
$d.Content.Find.Execute("50×95=", $true, $false, $false, $false, $false, $true, 1, $false, "57×61=", 2) | Out-Null
$d.Content.Find.Execute("77×69=", $true, $false, $false, $false, $false, $true, 1, $false, "54×76=", 2) | Out-Null
$d.Content.Find.Execute("97×43=", $true, $false, $false, $false, $false, $true, 1, $false, "64×24=", 2) | Out-Null
$d.Content.Find.Execute("68×19=", $true, $false, $false, $false, $false, $true, 1, $false, "77×20=", 2) | Out-Null
$d.Content.Find.Execute("53×35=", $true, $false, $false, $false, $false, $true, 1, $false, "22×41=", 2) | Out-Null

Write-Output "Replacements applied"